$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "61.762.45"
$ws.Range("E2").Value = "  -1.61%  "

# Row 3
$ws.Range("D3").Value = "2.893.04"
$ws.Range("E3").Value = "  -2.01%  "

# Row 4
$c = $ws.Range("D4")
$c.Value = "'1.00"
$c.Style = "Normal"
$ws.Range("E4").Value = "  +0.11%  "

# Row 5
$c = $ws.Range("D5")
$c.Value = "'568.11"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -4.35%  "

# Row 6
$c = $ws.Range("D6")
$c.Value = "'143.40"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -1.85%  "

# Row 7
$ws.Range("E7").Value = "  -0.07%  "

# Row 8
$c = $ws.Range("D8")
$c.Value = "'0.504"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -0.76%  "

# Row 9
$ws.Range("D9").Value = "2.890.60"
$ws.Range("E9").Value = "  -2.15%  "

# Row 10
$c = $ws.Range("D10")
$c.Value = "'6.62"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -8.88%  "

# Row 11
$ws.Range("E11").Value = "  -1.17%  "

# Row 12
$ws.Range("E12").Value = "  -2.66%  "

# Row 13
$c = $ws.Range("D13")
$c.Value = "'0.0000234"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -1.70%  "

# Row 14
$c = $ws.Range("D14")
$c.Value = "'31.94"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -3.41%  "

# Row 15
$ws.Range("E15").Value = "  -0.75%  "

# Row 16
$ws.Range("D16").Value = "3.372.02"
$ws.Range("E16").Value = "  -2.10%  "

# Row 17
$ws.Range("D17").Value = "61.716.82"
$ws.Range("E17").Value = "  -1.50%  "

# Row 18
$ws.Range("E18").Value = "  -1.88%  "

# Row 19
$ws.Range("D19").Value = "2.891.91"
$ws.Range("E19").Value = "  -2.78%  "

# Row 20
$c = $ws.Range("D20")
$c.Value = "'436.82"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -1.06%  "

# Row 21
$ws.Range("E21").Value = "  -2.39%  "

# Row 22
$c = $ws.Range("D22")
$c.Value = "'0.653"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -2.32%  "

# Row 23
$ws.Range("E23").Value = "  -2.67%  "

# Row 24
$c = $ws.Range("D24")
$c.Value = "'79.10"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -2.86%  "

# Row 25
$c = $ws.Range("D25")
$c.Value = "'11.91"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +0.25%  "

# Row 26
$c = $ws.Range("D26")
$c.Value = "'10.10"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -10.09%  "

# Row 28
$c = $ws.Range("D28")
$c.Value = "'2.03"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -4.65%  "

# Row 29
$ws.Range("E29").Value = "  +9.91%  "

# Row 30
$c = $ws.Range("D30")
$c.Value = "'7.05"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -3.08%  "

# Row 31
$c = $ws.Range("D31")
$c.Value = "'2.50"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -4.26%  "

# Row 32
$c = $ws.Range("D32")
$c.Value = "'2.06"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -4.45%  "

# Row 33
$c = $ws.Range("D33")
$c.Value = "'0.998"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -0.30%  "

# Row 34
$ws.Range("E34").Value = "  -2.82%  "

# Row 35
$c = $ws.Range("D35")
$c.Value = "'25.60"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -3.39%  "

# Row 36
$c = $ws.Range("D36")
$c.Value = "'0.953"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -3.73%  "

# Row 37
$ws.Range("E37").Value = "  -4.08%  "

# Row 38
$c = $ws.Range("D38")
$c.Value = "'48.98"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -1.14%  "

# Row 39
$c = $ws.Range("D39")
$c.Value = "'2.89"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -6.47%  "

# Row 40
$ws.Range("E40").Value = "  -4.68%  "

# Row 41
$ws.Range("E41").Value = "  -1.29%  "

# Row 42
$ws.Range("E42").Value = "  -3.27%  "

# Row 43
$ws.Range("B43").Value = "Arweave"
$ws.Range("C43").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$c = $ws.Range("D43")
$c.Value = "'39.24"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -4.22%  "

# Row 44
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$c = $ws.Range("D44")
$c.Value = "'0.268"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -4.42%  "

# Row 45
$ws.Range("D45").Value = "2.685.39"
$ws.Range("E45").Value = "  -1.52%  "

# Row 46
$c = $ws.Range("D46")
$c.Value = "'133.06"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -0.85%  "

# Row 47
$c = $ws.Range("D47")
$c.Value = "'0.0332"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -1.63%  "

# Row 48
$ws.Range("E48").Value = "  -0.03%  "

# Row 49
$c = $ws.Range("D49")
$c.Value = "'337.07"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -7.27%  "

# Row 50
$ws.Range("E50").Value = "  -2.11%  "

# Row 51
$c = $ws.Range("D51")
$c.Value = "'21.58"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -5.73%  "
